# "存款" (deposit) sheet: add bank/deposit_type/currency columns plus the
# standard property_category/category/date/legislator_name/legislator_id/
# source_file/index metadata columns (G:M), turn row 1 into real column
# headers (it used to just duplicate the first data row's values), and make
# column F ("total") hold a real number instead of a shared text string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# --- Row 1: column headers -------------------------------------------------
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# Column I holds a literal "2013-12-25" string (not an Excel date). Format as
# text first so the COM layer doesn't silently re-interpret it as a date
# serial, then restore the plain data-row style once the text is in place.
$ws.Range("I2:I8").NumberFormat = "@"

# --- Data rows 2-8 -----------------------------------------------------------
$rows = @(
    @{ Row=2; Index=78; Bank="集集郵局（第9支局）";         DepositType="活期存款";     Owner="陳麗珠"; Total=1609930 },
    @{ Row=3; Index=79; Bank="南投縣集集鎮農會信用部";       DepositType="活期存款";     Owner="陳麗珠"; Total=860288 },
    @{ Row=4; Index=80; Bank="南投縣集集鎮農會信用部";       DepositType="活期存款";     Owner="林明溱"; Total=354932 },
    @{ Row=5; Index=81; Bank="集集郵局（第9支局）";         DepositType="活期存款";     Owner="林明溱"; Total=19499 },
    @{ Row=6; Index=82; Bank="臺灣銀行南投分行";           DepositType="活期儲蓄存款"; Owner="林明溱"; Total=113369 },
    @{ Row=7; Index=83; Bank="合作金庫商業銀行集集分行";     DepositType="活期存款";     Owner="林明溱"; Total=726427 },
    @{ Row=8; Index=84; Bank="臺灣銀行南投分行";           DepositType="活期儲蓄存款"; Owner="陳麗珠"; Total=1079 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("B$row").Value = $r.Bank
    $ws.Range("C$row").Value = $r.DepositType
    $ws.Range("D$row").Value = "新臺幣"
    $ws.Range("E$row").Value = $r.Owner
    $ws.Range("F$row").Value = $r.Total
    $ws.Range("G$row").Value = "deposit"
    $ws.Range("H$row").Value = "normal"
    $ws.Range("I$row").Value = "2013-12-25"
    $ws.Range("J$row").Value = "林明溱"
    $ws.Range("K$row").Value = 1706
    $ws.Range("L$row").Value = "tmpfac21"
    $ws.Range("M$row").Value = $r.Index
}

# Put column I back on the same plain (unformatted) style as its neighbours
# now that the text is safely stored, so the number-format override doesn't
# linger as a visible difference from the rest of the row.
$ws.Range("I2:I8").Style = $ws.Range("H2").Style
